$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Header label for new column F (row 25), matching existing "Input in Deg" style header row.
$ws.Range("F25").Value = "Umdrehungen (Stange)"

# Fill formulas for F26:F44 -> E/(200*19.2), i.e. stepPosition starting at an
# initial value of 0 (0-based "Umdrehungen (Stange)" conversion of the step count).
for ($r = 26; $r -le 44; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Formula = "=E$r/(200*19.2)"
    $cell.NumberFormat = "General"
    $cell.HorizontalAlignment = -4131
    $cell.IndentLevel = 3
}

$ws.Range("I24").Select()
